$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    return $ok
}

function Append-EmptyParaAfterIndex($index) {
    $anchor = $d.Paragraphs.Item($index)
    $anchor.Range.InsertParagraphAfter()
}

function Append-TextParaAfterIndex($index, $text) {
    $anchor = $d.Paragraphs.Item($index)
    $anchor.Range.InsertParagraphAfter()
    $d.Paragraphs.Item($index + 1).Range.InsertBefore($text)
}

function Insert-TextParaBeforeIndex($index, $text) {
    $anchor = $d.Paragraphs.Item($index)
    $anchor.Range.InsertParagraphBefore()
    $d.Paragraphs.Item($index).Range.InsertBefore($text)
}

function Insert-EmptyParaBeforeIndex($index) {
    $anchor = $d.Paragraphs.Item($index)
    $anchor.Range.InsertParagraphBefore()
}

# ---------------------------------------------------------------------
# 1) Fix "сгенерированный" -> "сгенерированные" and add a new sentence
#    about having generated rules for 5000+ symptoms (paragraph 12,
#    in place - no paragraph-count change).
# ---------------------------------------------------------------------
Replace-Text `
    "библиотеку Spacy и заранее сгенерированный правила для выделения симптомов. Также" `
    "библиотеку Spacy и заранее сгенерированные правила для выделения симптомов. Мы сгенерировали правила для чуть больше 5000 симптомов. Также"

# ---------------------------------------------------------------------
# 2) Truncate paragraph 16 ("Ещё важно отметить...") - drop the trailing
#    sentence about training/validation data (in place).
# ---------------------------------------------------------------------
Replace-Text `
    " В качестве данных для обучения и валидации модели используются собранный датасет симптомов и болезней, которые выложен на платформу Kaggle, а также априорные знания о болезнях и их симптомов и слабоструктурированные данные с медицинских форумов." `
    ""

# ---------------------------------------------------------------------
# 3) Truncate paragraph 18 ("В результате разработана...") - drop the
#    trailing "На этом у меня всё, спасибо за внимание." sentence
#    (in place).
# ---------------------------------------------------------------------
Replace-Text `
    "На этом у меня всё, спасибо за внимание." `
    ""

# At this point the paragraph count is still 18 (no paragraphs were
# added/removed yet, only text inside existing runs changed).

# ---------------------------------------------------------------------
# 4) Append four new paragraphs after paragraph 18 (the now-truncated
#    "...на гитхабе. " paragraph), in order: empty, "future" paragraph,
#    empty, "thanks" paragraph.
# ---------------------------------------------------------------------
Append-EmptyParaAfterIndex 18
Append-TextParaAfterIndex 19 "В дальнейшем мы планируем замерить различные метрики качества инструмента. Поддержать работу с большим количеством симптомов. Улучшить точность распознавания симптомов. И провести эксперименты с более сложными моделями для определения симптомов."
Append-EmptyParaAfterIndex 20
Append-TextParaAfterIndex 21 "На этом у меня всё, спасибо за внимание."

# ---------------------------------------------------------------------
# 5) Insert two new paragraphs *before* paragraph 18 (the "В результате
#    разработана..." paragraph): the "metrics" paragraph, then an empty
#    paragraph, so they land right after the existing empty paragraph 17.
# ---------------------------------------------------------------------
Insert-EmptyParaBeforeIndex 18
Insert-TextParaBeforeIndex 18 "Также мы уже получили первые метрики качества нашего прототипа. Точность предсказания на тестовой выборке оказалась 0.23. Это довольно низкий показатель и мы уже работаем над его улучшением. Одно из направлений это повышать количество поддерживаемых симптомов и качество их распознавания. Другое направление это экспериментировать более сложные модели для предсказания диагноза по симптомам."

# ---------------------------------------------------------------------
# 6) Insert two new paragraphs *before* paragraph 10 ("На данном слайде
#    представлена архитектура модуля..."): the "dataset" paragraph, then
#    an empty paragraph.
# ---------------------------------------------------------------------
Insert-EmptyParaBeforeIndex 10
Insert-TextParaBeforeIndex 10 "В качестве Датасета мы используем априорные знания оболезнях и их симптомах, а также мы собрали свой датасет симптомов и болезней из форумов и других открытых источников. Этот датасет вложен в открытый доступ на платформе Kaggle и насчитывает чуть больше 5.000 записей. На диаграмме на слайде представлено распределение болезней по категориям. И например наибольшее количество болезней относится к категории дерматологии и венерологии."

Write-Output ("final paragraph count: " + $d.Paragraphs.Count)
